$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# B2 "ProcessABCQueue" -> "DermaItems" (same asset name cell, renamed in place)
$ws.Range("B2").Value = "DermaItems"

# B3: mark this config row as shared between processes
$ws.Range("B3").Value = "Shared"

# ---- New configuration rows 7-18 ----
$ws.Range("A7").Value  = "TemplateFilePath"
$ws.Range("A8").Value  = "AppendFilePath"
$ws.Range("A9").Value  = "AppendSheetName"
$ws.Range("B9").Value  = "Sheet1"
$ws.Range("A10").Value = "WebsiteExcelPath"
$ws.Range("A11").Value = "ProductExcelPath"
$ws.Range("A12").Value = "WebsiteSheetName"
$ws.Range("B12").Value = "Website"
$ws.Range("A13").Value = "ProductSheetName"
$ws.Range("B13").Value = "Product Details"
$ws.Range("A14").Value = "DermaUkSearch"
$ws.Range("A15").Value = "DermaFrSearch"
$ws.Range("B15").Value = "https://dermalogica.fr/pages/search-results-page?q="
$ws.Range("A16").Value = "LookFantasticSearch"
$ws.Range("B16").Value = "https://www.lookfantastic.com/elysium.search?search="
$ws.Range("A17").Value = "NewTabUrl"
$ws.Range("B17").Value = "chrome://settings/clearBrowserData"
$ws.Range("A18").Value = "ChromeNewTab"
$ws.Range("B18").Value = "chrome://newtab"

# UK search link added via Insert Hyperlink (populates the display text too)
$ws.Hyperlinks.Add($ws.Range("B14"), "https://www.dermalogica.co.uk/search?q=")

# File paths filled in last
$ws.Range("B7").Value  = "C:\Users\ShikharSaxena\Documents\UiPath\Dermalogica\Data\Input\TemplateProductDetails.xlsx"
$ws.Range("B8").Value  = "C:\Users\ShikharSaxena\Documents\UiPath\Dermalogica\Data\Output\TemplateProductDetails.xlsx"
$ws.Range("B10").Value = "C:\Users\ShikharSaxena\Documents\UiPath\Dermalogica\Data\Input\Websites.xlsx"
$ws.Range("B11").Value = "C:\Users\ShikharSaxena\Documents\UiPath\Dermalogica\Data\Input\Product Details.xlsx"

# FR / LookFantastic cells also get hyperlink relationships (text already matches)
$ws.Hyperlinks.Add($ws.Range("B15"), "https://dermalogica.fr/pages/search-results-page?q=")
$ws.Hyperlinks.Add($ws.Range("B16"), "https://www.lookfantastic.com/elysium.search?search=")

# Trim unused trailing rows (996-998) now that sheet content ends sooner
$ws.Rows("996:998").Delete()

# Make Settings the active sheet / restore the saved selection
$ws.Activate()
$ws.Range("B12").Select()
